# Update market-price derived columns (H-N) on several leve-profit rows
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets, per refreshed
# Universalis pricing data pulled by the scheduled runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 54: Arcane Arts for Dummies
$ws.Range("H54").Value = 10000
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 10000
$ws.Range("K54").Value = 0
$ws.Range("L54").Value = 10000
$ws.Range("M54").ClearContents()
$ws.Range("N54").Value = -10972
# Row 98: The Dotted Line
$ws.Range("H98").Value = 2058.5881
$ws.Range("I98").Value = 1676.1538
$ws.Range("J98").Value = 3301.5
$ws.Range("K98").Value = 1676.1538
$ws.Range("L98").Value = 3301.5
$ws.Range("M98").Value = -178.1538
$ws.Range("N98").Value = -6297.5
# Row 112: Making Ends Meet
$ws.Range("H112").Value = 3703
$ws.Range("I112").Value = 900
$ws.Range("J112").Value = 4904.2856
$ws.Range("K112").Value = 2700
$ws.Range("L112").Value = 14712.8568
$ws.Range("M112").Value = -1592
$ws.Range("N112").Value = -16928.8568
# Row 122: Wishful Inking
$ws.Range("H122").Value = 2058.5881
$ws.Range("I122").Value = 1676.1538
$ws.Range("J122").Value = 3301.5
$ws.Range("K122").Value = 5028.4614
$ws.Range("L122").Value = 9904.5
$ws.Range("M122").Value = -2578.4614
$ws.Range("N122").Value = -14804.5
# Row 138: All-night Crafting
$ws.Range("H138").Value = 1851.4651
$ws.Range("I138").Value = 1898.6538
$ws.Range("J138").Value = 1779.2941
$ws.Range("K138").Value = 5695.9614
$ws.Range("L138").Value = 5337.8823
$ws.Range("M138").Value = -555.9614000000001
$ws.Range("N138").Value = -15617.8823
# Row 141: Remedy for Reason
$ws.Range("H141").Value = 7015.6924
$ws.Range("I141").Value = 2179.2917
$ws.Range("J141").Value = 65052.5
$ws.Range("K141").Value = 6537.875100000001
$ws.Range("L141").Value = 195157.5
$ws.Range("M141").Value = -1357.875100000001
$ws.Range("N141").Value = -205517.5

$ws = $wb.Worksheets.Item("ARM")
# Row 32: Ingot We Trust
$ws.Range("H32").Value = 5304103
$ws.Range("I32").Value = 7436.1294
$ws.Range("J32").Value = 37084104
$ws.Range("K32").Value = 7436.1294
$ws.Range("L32").Value = 37084104
$ws.Range("M32").Value = -7149.1294
$ws.Range("N32").Value = -37084678
# Row 61: Dealing with the Tough Stuff
$ws.Range("H61").Value = 2236.4285
$ws.Range("I61").Value = 1196.24
$ws.Range("J61").Value = 4836.9
$ws.Range("K61").Value = 1196.24
$ws.Range("L61").Value = 4836.9
$ws.Range("M61").Value = -984.24
$ws.Range("N61").Value = -5260.9
# Row 132: Don't Bore Me, Ore Me
$ws.Range("H132").Value = 3576744
$ws.Range("I132").Value = 2531.4375
$ws.Range("J132").Value = 6586607.5
$ws.Range("K132").Value = 7594.3125
$ws.Range("L132").Value = 19759822.5
$ws.Range("M132").Value = -5064.3125
$ws.Range("N132").Value = -19764882.5
# Row 136: Metal with Mettle
$ws.Range("H136").Value = 2236.4285
$ws.Range("I136").Value = 1196.24
$ws.Range("J136").Value = 4836.9
$ws.Range("K136").Value = 3588.72
$ws.Range("L136").Value = 14510.7
$ws.Range("M136").Value = -1038.72
$ws.Range("N136").Value = -19610.7

$ws = $wb.Worksheets.Item("BSM")
# Row 94: High Steal
$ws.Range("H94").Value = 944.63635
$ws.Range("I94").Value = 944.63635
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 944.63635
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = -493.63635
$ws.Range("N94").ClearContents()
# Row 134: Ruthenium Supremium
$ws.Range("H134").Value = 5237
$ws.Range("I134").Value = 2205.92
$ws.Range("J134").Value = 9225.263000000001
$ws.Range("K134").Value = 6617.76
$ws.Range("L134").Value = 27675.789
$ws.Range("M134").Value = -4082.76
$ws.Range("N134").Value = -32745.789

$ws = $wb.Worksheets.Item("CRP")
# Row 16: Raise the Roof
$ws.Range("H16").Value = 2661.2354
$ws.Range("I16").Value = 3212.8
$ws.Range("J16").Value = 1873.2858
$ws.Range("K16").Value = 3212.8
$ws.Range("L16").Value = 1873.2858
$ws.Range("M16").Value = -2925.8
$ws.Range("N16").Value = -2447.2858
# Row 31: Wall Not Found
$ws.Range("H31").Value = 7144766
$ws.Range("I31").Value = 1351.625
$ws.Range("J31").Value = 16669319
$ws.Range("K31").Value = 1351.625
$ws.Range("L31").Value = 16669319
$ws.Range("M31").Value = -1056.625
$ws.Range("N31").Value = -16669909
# Row 34: Armoires of the Rich and Famous
$ws.Range("H34").Value = 7144766
$ws.Range("I34").Value = 1351.625
$ws.Range("J34").Value = 16669319
$ws.Range("K34").Value = 1351.625
$ws.Range("L34").Value = 16669319
$ws.Range("M34").Value = -1149.625
$ws.Range("N34").Value = -16669723
# Row 58: You Do the Heavy Lifting
$ws.Range("H58").Value = 3346759.2
$ws.Range("I58").Value = 5798.143
$ws.Range("J58").Value = 6270100
$ws.Range("K58").Value = 5798.143
$ws.Range("L58").Value = 6270100
$ws.Range("M58").Value = -5595.143
$ws.Range("N58").Value = -6270506
# Row 81: Don't Ask Wyvern
$ws.Range("H81").Value = 68800
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 68800
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 68800
$ws.Range("N81").Value = -70796
# Row 84: A Sky Pirate's Life for Me (L)
$ws.Range("H84").Value = 68800
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 68800
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 206400
$ws.Range("N84").Value = -216384
# Row 113: Patient Patients
$ws.Range("H113").Value = 2661.2354
$ws.Range("I113").Value = 3212.8
$ws.Range("J113").Value = 1873.2858
$ws.Range("K113").Value = 3212.8
$ws.Range("L113").Value = 1873.2858
$ws.Range("M113").Value = -1042.8
$ws.Range("N113").Value = -6213.2858
# Row 136: Turali Quality
$ws.Range("H136").Value = 3346759.2
$ws.Range("I136").Value = 5798.143
$ws.Range("J136").Value = 6270100
$ws.Range("K136").Value = 17394.429
$ws.Range("L136").Value = 18810300
$ws.Range("M136").Value = -14844.429
$ws.Range("N136").Value = -18815400
# Row 141: No Greater Treasure
$ws.Range("H141").Value = 39994
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 39994
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 39994
$ws.Range("N141").Value = -50354

$ws = $wb.Worksheets.Item("CUL")
# Row 107: Slippery Service
$ws.Range("H107").Value = 11346155
$ws.Range("I107").Value = 26317262
$ws.Range("J107").Value = 405729.53
$ws.Range("K107").Value = 78951786
$ws.Range("L107").Value = 1217188.59
$ws.Range("M107").Value = -78949866
$ws.Range("N107").Value = -1221028.59
# Row 136: Simple Is Hardest
$ws.Range("H136").Value = 3666.682
$ws.Range("I136").Value = 2888
$ws.Range("J136").Value = 4445.364
$ws.Range("K136").Value = 8664
$ws.Range("L136").Value = 13336.092
$ws.Range("M136").Value = -3564
$ws.Range("N136").Value = -23536.092

$ws = $wb.Worksheets.Item("GSM")
# Row 97: If I'd a Koppranickel for Every Time...
$ws.Range("H97").Value = 1386.8214
$ws.Range("I97").Value = 1453.3334
$ws.Range("J97").Value = 1187.2858
$ws.Range("K97").Value = 1453.3334
$ws.Range("L97").Value = 1187.2858
$ws.Range("M97").Value = -957.3334
$ws.Range("N97").Value = -2179.2858

$ws = $wb.Worksheets.Item("LTW")
# Row 93: Hide to Go Seek
$ws.Range("H93").Value = 279395.78
$ws.Range("I93").Value = 358597.84
$ws.Range("J93").Value = 2188.5
$ws.Range("K93").Value = 358597.84
$ws.Range("L93").Value = 2188.5
$ws.Range("M93").Value = -357349.84
$ws.Range("N93").Value = -4684.5
# Row 136: Respect for Br'aax
$ws.Range("H136").Value = 33334964
$ws.Range("I136").Value = 50001172
$ws.Range("J136").Value = 2554
$ws.Range("K136").Value = 150003516
$ws.Range("L136").Value = 7662
$ws.Range("M136").Value = -150000966
$ws.Range("N136").Value = -12762

$ws = $wb.Worksheets.Item("WVR")
# Row 122: Heavy Armoire
$ws.Range("H122").Value = 8090.2085
$ws.Range("I122").Value = 9429.412
$ws.Range("J122").Value = 4837.857
$ws.Range("K122").Value = 28288.236
$ws.Range("L122").Value = 14513.571
$ws.Range("M122").Value = -25838.236
$ws.Range("N122").Value = -19413.571
